$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "return to home page" troubleshooting steps to column L
$ws.Range("L19").Value = "when wont let u use flask url"
$ws.Range("L20").Value = "go to this link chrome://net-internals/#hsts"
$ws.Range("L21").Value = "go to delete domain security policies"
$ws.Range("L22").Value = "type in 127.0.0.1"
$ws.Range("L23").Value = "delete and reload link"

# Move the active selection down to the next empty row in column L
$ws.Range("L24").Select()
